# Atualizações e correções de artefatos
# Adds a new "Fase de Elaboração E1" block (banner + header + divider + 9 task
# rows) below the existing "Lista de Itens de Trabalho" table.

function RGBColor($r, $g, $b) { return $r + ($g * 256) + ($b * 65536) }

$COLOR_YELLOW = RGBColor 0xFF 0xFF 0x00
$COLOR_GRAY_FG = RGBColor 0xA6 0xA6 0xA6
$COLOR_GRAY_BG = RGBColor 0xC0 0xC0 0xC0
$COLOR_CREAM_FG = RGBColor 0xFF 0xFF 0xFF
$COLOR_CREAM_BG = RGBColor 0xFF 0xFF 0xCC

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lista de Itens de Trabalho")

# Restore default (visible) gridlines on the active window/sheet.
$excel.ActiveWindow.DisplayGridlines = $true

# ---------------------------------------------------------------------------
# Row 15:16 - merged yellow banner
# ---------------------------------------------------------------------------
$banner = $ws.Range("A15:I16")
$banner.Interior.Color = $COLOR_YELLOW
$banner.Interior.PatternColor = $COLOR_YELLOW
$banner.Font.Bold = $true
$banner.HorizontalAlignment = -4108
$banner.VerticalAlignment = -4108
$ws.Rows(15).RowHeight = 12.85
$ws.Rows(16).RowHeight = 12.85
$ws.Range("A15:I16").Merge()
$ws.Range("A15").Value = "Fase de Elaboração E1"

# ---------------------------------------------------------------------------
# Row 17 - repeated table header
# ---------------------------------------------------------------------------
$headerVals = @("Nome / Descrição", "Prioridade", "Tamanho Estimado (Pontos)", "Estado Atual", "Iteração Alvo", "Atribuído a", "Esforço Estimado (horas)", "Horas Trabalhadas", "Material de Referência")
for ($i = 0; $i -lt 9; $i++) {
    $ws.Cells.Item(17, $i + 1).Value = $headerVals[$i]
}
$headerRng = $ws.Range("A17:I17")
$headerRng.Interior.Color = $COLOR_GRAY_FG
$headerRng.Interior.PatternColor = $COLOR_GRAY_BG
$headerRng.Font.Bold = $true
$headerRng.Borders.LineStyle = 1
$headerRng.HorizontalAlignment = -4108
$ws.Rows(17).RowHeight = 12.85

# ---------------------------------------------------------------------------
# Row 18 - cream divider row
# ---------------------------------------------------------------------------
$divider = $ws.Range("A18:I18")
$divider.Interior.Color = $COLOR_CREAM_FG
$divider.Interior.PatternColor = $COLOR_CREAM_BG
$divider.HorizontalAlignment = -4108
for ($col = 1; $col -le 9; $col++) {
    $ws.Cells.Item(18, $col).Value = ""
}
$ws.Range("C18").Font.Bold = $true
$ws.Range("G18").Font.Bold = $true
$ws.Range("H18").Font.Bold = $true
$ws.Rows(18).RowHeight = 12.85

# ---------------------------------------------------------------------------
# Rows 19:27 - task data
# ---------------------------------------------------------------------------
# Columns: A Nome/Descrição, B Prioridade, C Tamanho Estimado, D Estado Atual,
#          E Iteração Alvo, F Atribuído a, G Esforço Estimado, H Horas
#          Trabalhadas, I Material de Referência (all blank here).
$taskRows = @(
    @("Alteração Plano de Projeto  ", "Alta", 2, "Iniciado", "E1", "Waltson", 2, 2),
    @("Criação do  Plano de Interação  E1", "Alta", 1, "Iniciado", "E1", "Waltson", 2, 2),
    @("Atualizar Lista de Riscos ", "Alta", 1, "Iniciado", "E1", "Ivson", "1/2 de hora", 1),
    @("Alteração Lista de Itens", "Alta", 1, "Iniciado", "E1", "Ivson", 1, 1),
    @("Refinamento Requisitos de Cadastro", "Alta", 1, "Iniciado", "E1", "waltson/Ivson", 2, 1),
    @("Refinamento Requisito de Entrega ", "Alta", 2, "Iniciado", "E1", "waltson/Ivson", 2, 2),
    @("Implementação dos Requisitos de Cadastro", "Alta", 6, "Em andamento", "E1", "waltson/Ivson", 5, 2),
    @("Implementação dos Requisitos de Entrega", "Alta", 5, "Em andamento", "E1", "waltson/Ivson", 8, 1),
    @("Modelagem de analise e projeto", "Alta", 5, "Em andamento", "E1", "waltson/Ivson", 8, 1)
)

$startRow = 19
for ($i = 0; $i -lt $taskRows.Count; $i++) {
    $r = $startRow + $i
    $data = $taskRows[$i]
    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = $data[4]
    $ws.Cells.Item($r, 6).Value = $data[5]
    $ws.Cells.Item($r, 7).Value = $data[6]
    $ws.Cells.Item($r, 8).Value = $data[7]

    $ws.Range("B$r:H$r").HorizontalAlignment = -4108

    if ($r -eq 23) {
        # Matches a minor inconsistency present in the authored workbook:
        # this single cell keeps the default (general/bottom) alignment.
        $ws.Cells.Item($r, 9).HorizontalAlignment = -4142
        $ws.Cells.Item($r, 9).VerticalAlignment = -4107
    }
    else {
        $ws.Cells.Item($r, 9).HorizontalAlignment = -4108
        $ws.Cells.Item($r, 9).VerticalAlignment = -4108
    }

    $ws.Rows($r).RowHeight = 12.85
}

# Column I of the task rows has no content, only formatting - touch it so the
# cell exists in the sheet even though it stays empty.
for ($i = 0; $i -lt $taskRows.Count; $i++) {
    $r = $startRow + $i
    if ($ws.Cells.Item($r, 9).Value -eq $null) {
        $ws.Cells.Item($r, 9).NumberFormat = "General"
    }
}

$ws.Range("I33").Select()
Write-Output "done"
